$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each touched row, force Text number format on the touched
# columns before writing values, so numeric-looking strings (e.g.
# "95.40", "1.79") are preserved verbatim as text instead of being
# auto-converted to numbers by Excel. ClearFormats() afterward removes
# the temporary format so styling matches the original (no explicit s=).

$rowRange = $ws.Range('D2:E2')
$rowRange.NumberFormat = "@"
$ws.Range('D2').Value = '42.801.21'
$ws.Range('E2').Value = '  -0.08%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D3:E3')
$rowRange.NumberFormat = "@"
$ws.Range('D3').Value = '2.316.39'
$ws.Range('E3').Value = '  +0.51%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E4')
$rowRange.NumberFormat = "@"
$ws.Range('E4').Value = '  -0.02%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E5')
$rowRange.NumberFormat = "@"
$ws.Range('E5').Value = '  -1.23%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D6:E6')
$rowRange.NumberFormat = "@"
$ws.Range('D6').Value = '95.40'
$ws.Range('E6').Value = '  -1.48%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E7')
$rowRange.NumberFormat = "@"
$ws.Range('E7').Value = '  -0.17%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E8')
$rowRange.NumberFormat = "@"
$ws.Range('E8').Value = '  -0.01%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E9')
$rowRange.NumberFormat = "@"
$ws.Range('E9').Value = '  -1.27%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D10:E10')
$rowRange.NumberFormat = "@"
$ws.Range('D10').Value = '34.10'
$ws.Range('E10').Value = '  -3.22%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D11:E11')
$rowRange.NumberFormat = "@"
$ws.Range('D11').Value = '18.92'
$ws.Range('E11').Value = '  +1.59%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D12:E12')
$rowRange.NumberFormat = "@"
$ws.Range('D12').Value = '0.0783'
$ws.Range('E12').Value = '  -0.28%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E13')
$rowRange.NumberFormat = "@"
$ws.Range('E13').Value = '  +0.40%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E14')
$rowRange.NumberFormat = "@"
$ws.Range('E14').Value = '  -2.02%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D15:E15')
$rowRange.NumberFormat = "@"
$ws.Range('D15').Value = '2.674.83'
$ws.Range('E15').Value = '  +0.46%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D16:E16')
$rowRange.NumberFormat = "@"
$ws.Range('D16').Value = '2.348.14'
$ws.Range('E16').Value = '  +2.33%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D17')
$rowRange.NumberFormat = "@"
$ws.Range('D17').Value = '0.787'
$rowRange.ClearFormats()

$rowRange = $ws.Range('D18')
$rowRange.NumberFormat = "@"
$ws.Range('D18').Value = '42.739.80'
$rowRange.ClearFormats()

$rowRange = $ws.Range('E19')
$rowRange.NumberFormat = "@"
$ws.Range('E19').Value = '  -4.42%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D20:E20')
$rowRange.NumberFormat = "@"
$ws.Range('D20').Value = '6.14'
$ws.Range('E20').Value = '  +1.89%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D21:E21')
$rowRange.NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0889'
$ws.Range('E21').Value = '  -0.53%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D22:E22')
$rowRange.NumberFormat = "@"
$ws.Range('D22').Value = '67.73'
$ws.Range('E22').Value = '  +0.80%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E23')
$rowRange.NumberFormat = "@"
$ws.Range('E23').Value = '  +6.36%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D24:E24')
$rowRange.NumberFormat = "@"
$ws.Range('D24').Value = '235.44'
$ws.Range('E24').Value = '  -0.22%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E25')
$rowRange.NumberFormat = "@"
$ws.Range('E25').Value = '  -0.13%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D26:E26')
$rowRange.NumberFormat = "@"
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  +0.54%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D27:E27')
$rowRange.NumberFormat = "@"
$ws.Range('D27').Value = '24.36'
$ws.Range('E27').Value = '  -1.42%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E28')
$rowRange.NumberFormat = "@"
$ws.Range('E28').Value = '  +14.64%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E29')
$rowRange.NumberFormat = "@"
$ws.Range('E29').Value = '  +0.80%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D30:E30')
$rowRange.NumberFormat = "@"
$ws.Range('D30').Value = '32.15'
$ws.Range('E30').Value = '  -2.93%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D31:E31')
$rowRange.NumberFormat = "@"
$ws.Range('D31').Value = '147.92'
$ws.Range('E31').Value = '  -11.02%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E32')
$rowRange.NumberFormat = "@"
$ws.Range('E32').Value = '  -0.08%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E33')
$rowRange.NumberFormat = "@"
$ws.Range('E33').Value = '  +0.31%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D34:E34')
$rowRange.NumberFormat = "@"
$ws.Range('D34').Value = '17.65'
$ws.Range('E34').Value = '  -2.11%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E35')
$rowRange.NumberFormat = "@"
$ws.Range('E35').Value = '  -0.11%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E36')
$rowRange.NumberFormat = "@"
$ws.Range('E36').Value = '  +2.03%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E37')
$rowRange.NumberFormat = "@"
$ws.Range('E37').Value = '  -1.22%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D38:E38')
$rowRange.NumberFormat = "@"
$ws.Range('D38').Value = '1.79'
$ws.Range('E38').Value = '  +2.71%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E39')
$rowRange.NumberFormat = "@"
$ws.Range('E39').Value = '  -0.58%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D40:E40')
$rowRange.NumberFormat = "@"
$ws.Range('D40').Value = '2.73'
$ws.Range('E40').Value = '  +0.62%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('B41:E41')
$rowRange.NumberFormat = "@"
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = '22.06'
$ws.Range('E41').Value = '  +21.75%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('B42:E42')
$rowRange.NumberFormat = "@"
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = '0.108'
$ws.Range('E42').Value = '  -0.95%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D43:E43')
$rowRange.NumberFormat = "@"
$ws.Range('D43').Value = '1.917.17'
$ws.Range('E43').Value = '  -4.08%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E44')
$rowRange.NumberFormat = "@"
$ws.Range('E44').Value = '  -0.89%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D45:E45')
$rowRange.NumberFormat = "@"
$ws.Range('D45').Value = '10.07'
$ws.Range('E45').Value = '  -1.50%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E47')
$rowRange.NumberFormat = "@"
$ws.Range('E47').Value = '  -1.08%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('E48')
$rowRange.NumberFormat = "@"
$ws.Range('E48').Value = '  +1.51%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D49:E49')
$rowRange.NumberFormat = "@"
$ws.Range('D49').Value = '2.543.59'
$ws.Range('E49').Value = '  +0.51%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D50:E50')
$rowRange.NumberFormat = "@"
$ws.Range('D50').Value = '53.27'
$ws.Range('E50').Value = '  -0.52%  '
$rowRange.ClearFormats()

$rowRange = $ws.Range('D51:E51')
$rowRange.NumberFormat = "@"
$ws.Range('D51').Value = '72.22'
$ws.Range('E51').Value = '  +1.56%  '
$rowRange.ClearFormats()
